# Generate Report for Handoff
# Update the localization-status workbook to reflect that e2e\b.md has been
# handed off for localization: status moves from "Handed back: in sync with
# en-US" to "Ready for handoff", new handoff xliff files are recorded, and an
# error detail note is attached because the handback version is stale.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$handoffDate = "2016-08-18 14:40:26"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6aa3b8f9351472a280ceec243cbcbdbea5dc6f1/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1184c2631e16ffad830336edcc54b6a78ad1a5c/e2e/b.md."

# --- Overview sheet: row 3 is e2e\b.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = $handoffDate

# --- zh-cn sheet: row 3 is b.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-18 14:40:19"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").ColumnWidth = 39.17

# --- de-de sheet: row 3 is b.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $handoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.17
